# Update PureBatsmen_ODI sheet: replace India (IND) data with Australia (AUS) data.
# The new data set only has 9 values (rows 2-10) instead of the previous 14
# (rows 2-15), so the trailing rows are cleared to shrink the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country code header
$ws.Range("A1").Value = "AUS"

# Numeric-looking values need to stay text, like the original sheet, so a
# leading apostrophe is used to force text storage instead of numbers.
$ws.Range("A2").Value  = "'4369"
$ws.Range("A3").Value  = "'4726"
$ws.Range("A4").Value  = "'4558"
$ws.Range("A5").Value  = "'6471"
$ws.Range("A6").Value  = "'3842"
$ws.Range("A7").Value  = "'4824"
$ws.Range("A8").Value  = "'3725"
$ws.Range("A9").Value  = "'5860"
$ws.Range("A10").Value = "'3910"

# Remove the now-unused trailing rows (11-15) so the sheet's used range /
# dimension shrinks back down to A1:A10.
$ws.Range("A11:A15").ClearContents() | Out-Null
